# Solve Leetcode - 994. Rotting Oranges - Multi Source BFS
# Adds a new row (34) to the "Neetcode 150" sheet describing the newly
# solved problem, matching the existing table's layout/styling.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Neetcode 150")

$category = "Graphs"
$difficulty = "Medium"
$title = "994. Rotting Oranges"
$notes = "Just run a multi source bfs, and track the number of fresh fruits. while q and fruits > 0, then go down level by level in bfs with for _ in range(len(q)) and increment time everytime we go down a level. In the end return -1 if there are fresh fruits left, otherwise return time"

$row = 34

$ws.Cells.Item($row, 1).Value = $category
$ws.Cells.Item($row, 2).Value = $difficulty
$ws.Cells.Item($row, 3).Value = $title
$ws.Cells.Item($row, 4).Value = $notes

# Mirror the style of the row directly above (row 33) for each column.
$ws.Cells.Item($row, 2).Style = $ws.Cells.Item($row - 1, 2).Style
$ws.Cells.Item($row, 3).Style = $ws.Cells.Item($row - 1, 3).Style
$ws.Cells.Item($row, 4).Style = $ws.Cells.Item($row - 1, 4).Style

$ws.Rows.Item($row).RowHeight = $ws.Rows.Item($row - 1).RowHeight

$ws.Hyperlinks.Add(
    $ws.Cells.Item($row, 3),
    "https://leetcode.com/problems/rotting-oranges/",
    "",
    "",
    $title
)

$ws.Range("D35").Select()
